$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 95, pushing existing rows 95-156 down to 97-158.
$ws.Rows.Item(95).Insert()
$ws.Rows.Item(95).Insert()

# Populate new row 95 (Black Amber / Primera)
$ws.Range("A95").Value = 10
$ws.Range("B95").Value = "Vega Modelo de Temuco"
$ws.Range("C95").Value = "La Araucanía"
$ws.Range("D95").Value = 44574
$ws.Range("E95").Value = 9
$ws.Range("F95").Value = "Fruta"
$ws.Range("G95").Value = 100103
$ws.Range("H95").Value = "Frutos de hueso (carozo)"
$ws.Range("I95").Value = 100103002
$ws.Range("J95").Value = "Ciruela"
$ws.Range("K95").Value = "Black Amber"
$ws.Range("L95").Value = "Primera"
$ws.Range("M95").Value = 500
$ws.Range("N95").Value = 16000
$ws.Range("O95").Value = 16000
$ws.Range("P95").Value = 16000
$ws.Range("Q95").Value = "$/bandeja 18 kilos granel"
$ws.Range("R95").Value = "Región de O'Higgins"
$ws.Range("S95").Value = 889
$ws.Range("T95").Value = 18

# Populate new row 96 (Black Amber / Segunda)
$ws.Range("A96").Value = 10
$ws.Range("B96").Value = "Vega Modelo de Temuco"
$ws.Range("C96").Value = "La Araucanía"
$ws.Range("D96").Value = 44574
$ws.Range("E96").Value = 9
$ws.Range("F96").Value = "Fruta"
$ws.Range("G96").Value = 100103
$ws.Range("H96").Value = "Frutos de hueso (carozo)"
$ws.Range("I96").Value = 100103002
$ws.Range("J96").Value = "Ciruela"
$ws.Range("K96").Value = "Black Amber"
$ws.Range("L96").Value = "Segunda"
$ws.Range("M96").Value = 100
$ws.Range("N96").Value = 12000
$ws.Range("O96").Value = 12000
$ws.Range("P96").Value = 12000
$ws.Range("Q96").Value = "$/bandeja 18 kilos granel"
$ws.Range("R96").Value = "Región del Maule"
$ws.Range("S96").Value = 667
$ws.Range("T96").Value = 18

Write-Host "done"
